# Fix exploration time formula (#13)
# Recomputed Avg_Agent_Step_Time (G), Avg_Experiment_Time (H),
# Std_Agent_Step_Time (M) and Std_Experiment_Time (N) for rows 2-13
# on the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = @{ G = 2.43753139;          H = 166.85076695;         M = 0.2808116673432284;  N = 36.96750933532022 }
    3  = @{ G = 1.43362922;          H = 159.68834131;         M = 0.3857751479011183;  N = 50.99682913882356 }
    4  = @{ G = 1.26306222;          H = 56.55318427;          M = 0.2686290274461385;  N = 21.46195678839939 }
    5  = @{ G = 0.730662;            H = 47.42466252;          M = 0.1994292735205565;  N = 18.003032351658 }
    6  = @{ G = 0.57986558;          H = 16.04977179;          M = 0.1691694069141915;  N = 7.855137124474721 }
    7  = @{ G = 0.34253878;          H = 11.79694189;          M = 0.09894841761271202; N = 4.939485261636971 }
    8  = @{ G = 0.33216157;          H = 6.24381025;           M = 0.1163972933660317;  N = 3.833098956267941 }
    9  = @{ G = 0.19160893;          H = 4.31384833;           M = 0.05772261005937939; N = 1.920355393550123 }
    10 = @{ G = 0.20314113;          H = 2.75909081;           M = 0.08173785510263822; N = 1.904877334205648 }
    11 = @{ G = 0.12062729;          H = 2.02756292;           M = 0.04237322676056193; N = 1.04764165607101 }
    12 = @{ G = 0.12852231;          H = 1.31716232;           M = 0.0569629458522298;  N = 0.9971903187449306 }
    13 = @{ G = 0.08417903999999998; H = 1.14681178;           M = 0.03254535661374008; N = 0.6719166098999868 }
}

foreach ($row in $values.Keys) {
    $rowVals = $values[$row]
    $ws.Range("G$row").Value = $rowVals.G
    $ws.Range("H$row").Value = $rowVals.H
    $ws.Range("M$row").Value = $rowVals.M
    $ws.Range("N$row").Value = $rowVals.N
}
